$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NFT- RA")

# Row 4 - update assessment values
$ws.Range("E4").Value = "Moderate"
$ws.Range("G4").Value = "Moderate"
$ws.Range("H4").Value = ">30 to 50 %"
$ws.Range("I4").Value = "GREEN"

# Row 4 height
$ws.Rows.Item(4).RowHeight = 120

# Approvals / Signoff updates (written before J4 so shared-string order matches)
$ws.Range("H12").Value = "Approved by Preethi Govindaraj"
$ws.Range("I16").Value = "Approved by Preethi Govindaraj"

$ws.Range("J4").Value = "The design of the monitoring system is robust. Data throughput is handled well by the system. As per deployment of Mobile Application"

# Update selection to G4
$ws.Range("G4").Select()

$wb.Save()
